$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.164
$ws.Range("C2").Value = 7.891
$ws.Range("D2").Value = 29.211
$ws.Range("E2").Value = 11.225
$ws.Range("F2").Value = 6.694
$ws.Range("G2").Value = 17.726
$ws.Range("H2").Value = -5.06
$ws.Range("I2").Value = -1.291
$ws.Range("J2").Value = 11.976
$ws.Range("K2").Value = 11.131
$ws.Range("L2").Value = 5.575
$ws.Range("B3").Value = -10.333
$ws.Range("C3").Value = -9.803000000000001
$ws.Range("D3").Value = -20.32
$ws.Range("E3").Value = -12.271
$ws.Range("F3").Value = -20.33
$ws.Range("G3").Value = -6.338
$ws.Range("H3").Value = -40.304
$ws.Range("I3").Value = -25.764
$ws.Range("J3").Value = -1.991
$ws.Range("K3").Value = -9.882999999999999
$ws.Range("L3").Value = -16.511
$ws.Range("B4").Value = -11.331
$ws.Range("C4").Value = -10.105
$ws.Range("D4").Value = -24.831
$ws.Range("E4").Value = -8.999000000000001
$ws.Range("F4").Value = -16.973
$ws.Range("G4").Value = -5.91
$ws.Range("H4").Value = -37.973
$ws.Range("I4").Value = -19.977
$ws.Range("K4").Value = -4.593
$ws.Range("L4").Value = -14.185
$ws.Range("B5").Value = -15.759
$ws.Range("C5").Value = -3.63
$ws.Range("D5").Value = -30.538
$ws.Range("E5").Value = 6.471
$ws.Range("F5").Value = -13.576
$ws.Range("G5").Value = 4.67
$ws.Range("H5").Value = -31.513
$ws.Range("I5").Value = -7.868
$ws.Range("J5").Value = 2.705
$ws.Range("K5").Value = 4.363
$ws.Range("L5").Value = -20.876
$ws.Range("B6").Value = -9.382999999999999
$ws.Range("C6").Value = 44.107
$ws.Range("D6").Value = -31.303
$ws.Range("E6").Value = 52.502
$ws.Range("F6").Value = 34.755
$ws.Range("G6").Value = 21.981
$ws.Range("H6").Value = -34.861
$ws.Range("I6").Value = 2.625
$ws.Range("K6").Value = 74.58499999999999
$ws.Range("L6").Value = -20.003
$ws.Range("B7").Value = -2.129
$ws.Range("C7").Value = 8.308
$ws.Range("D7").Value = -7.874
$ws.Range("E7").Value = 9.656000000000001
$ws.Range("F7").Value = 6.733
$ws.Range("G7").Value = 4.436
$ws.Range("H7").Value = -8.938000000000001
$ws.Range("I7").Value = 0.5679999999999999
$ws.Range("J7").Value = 3.751
$ws.Range("K7").Value = 12.944
$ws.Range("L7").Value = -4.758
$ws.Range("B8").Value = 96201
$ws.Range("C8").Value = 178770.71
$ws.Range("D8").Value = 68542.36
$ws.Range("E8").Value = 71361.67999999999
$ws.Range("F8").Value = 157342.24
$ws.Range("G8").Value = 139564.24
$ws.Range("H8").Value = 77182.3
$ws.Range("I8").Value = 60216.84
$ws.Range("K8").Value = 215593.34
$ws.Range("L8").Value = 76895.56
$ws.Range("D9").Value = 1.045
$ws.Range("F9").Value = 1.007
$ws.Range("G9").Value = 1.032
$ws.Range("H9").Value = 1.115
$ws.Range("J9").Value = 0.615
$ws.Range("K9").Value = 1.034
$ws.Range("L9").Value = 0.608
$ws.Range("B10").Value = -9.423
$ws.Range("C10").Value = 0.428
$ws.Range("D10").Value = -16.206
$ws.Range("E10").Value = 1.842
$ws.Range("F10").Value = -1.322
$ws.Range("G10").Value = -3.803
$ws.Range("H10").Value = -17.786
$ws.Range("I10").Value = -6.402
$ws.Range("J10").Value = -1.411
$ws.Range("K10").Value = 4.691
$ws.Range("L10").Value = -9.863
$ws.Range("B11").Value = -0.135
$ws.Range("C11").Value = 0.381
$ws.Range("D11").Value = -0.263
$ws.Range("E11").Value = 0.371
$ws.Range("F11").Value = 0.278
$ws.Range("G11").Value = 0.157
$ws.Range("H11").Value = -0.35
$ws.Range("I11").Value = -0.002
$ws.Range("J11").Value = 0.221
$ws.Range("K11").Value = 0.519
$ws.Range("L11").Value = -0.244
$ws.Range("B12").Value = -0.03
$ws.Range("C12").Value = 0.078
$ws.Range("D12").Value = -0.081
$ws.Range("E12").Value = 0.093
$ws.Range("F12").Value = 0.061
$ws.Range("G12").Value = 0.037
$ws.Range("H12").Value = -0.08599999999999999
$ws.Range("I12").Value = -0.001
$ws.Range("J12").Value = 0.051
$ws.Range("K12").Value = 0.119
$ws.Range("L12").Value = -0.08799999999999999
$ws.Range("B14").Value = 20.319
$ws.Range("C14").Value = 20.18
$ws.Range("D14").Value = 32.284
$ws.Range("E14").Value = 24.364
$ws.Range("F14").Value = 21.99
$ws.Range("G14").Value = 24.366
$ws.Range("H14").Value = 27.348
$ws.Range("I14").Value = 24.676
$ws.Range("J14").Value = 14.176
$ws.Range("K14").Value = 23.761
$ws.Range("L14").Value = 22.093
$ws.Range("B15").Value = 0.421
$ws.Range("E15").Value = 0.407
$ws.Range("F15").Value = 0.785
$ws.Range("H15").Value = 0.464
$ws.Range("I15").Value = 0.316
$ws.Range("J15").Value = 0.575
$ws.Range("K15").Value = 0.741
$ws.Range("L15").Value = 0.863
$ws.Range("B16").Value = 7.294
$ws.Range("C16").Value = 7.88
$ws.Range("D16").Value = 8.332000000000001
$ws.Range("F16").Value = 8.054
$ws.Range("G16").Value = 8.239000000000001
$ws.Range("H16").Value = 8.848000000000001
$ws.Range("I16").Value = 6.969
$ws.Range("J16").Value = 5.161
$ws.Range("K16").Value = 8.253
$ws.Range("L16").Value = 5.105
